# Actualización automatica mar abr  6 17:33:39 CEST 2021
# The column "porcentaje-participacion" / "iaest-measure:porcentaje-participacion"
# metadata row entries are replaced by duplicates of the "participacion" /
# "iaest-measure:participacion" values (column L now mirrors column J), and the
# now-unused shared strings are dropped from the workbook when Excel re-saves it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = $ws.Range("J2").Value2
$ws.Range("L3").Value = $ws.Range("J3").Value2
